# The authors' "Fehler" category previously used a value with literal
# "\n" escape sequences as line separators ("Kompromittierung\nvon\n
# Geschäfts-E-Mails"). The edit normalizes this label to a single-line,
# space-separated string: "Kompromittierung von Geschäfts-E-Mails".
# (Internally Excel re-sorts the shared-strings table on save, which is
# why the raw XML diff also shows a lot of <si> index churn elsewhere —
# but the only actual content change is this one cell's text.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Kompromittierung von Geschäfts-E-Mails"
